$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PCB designation text (Vier_Gewinnt V1.0 -> Vier_Gewinnt V2.1)
$ws.Range("C4").Value = "Vier_Gewinnt V2.1"

# Remove the three obsolete BOM rows (delete entire rows, shifting cells up).
# Order: delete from bottom to top so earlier row numbers stay valid.
# Row 18: Elektrolytkondensator (1 µF)
$ws.Rows.Item(18).Delete()
# Row 15: Diode (1N4007)
$ws.Rows.Item(15).Delete()
# Row 6: 5er 221-WAGO-Klemmen (6mm2)
$ws.Rows.Item(6).Delete()

# Update the selected range / active cell to match the new table extent
$ws.Range("A3:C15").Select()
